# top_tax_affected_positive_various.xlsx
# Add a "Russia" country column (inserted before Saudi Arabia, i.e. new column L,
# pushing Saudi Arabia -> M and USA -> N), refresh the whole data table with the
# new survey figures, and fix the "Affected by ..." (any variant) label text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert "Russia" before "Saudi Arabia", "USA" moves to new col N ---
$ws.Range("L1").Value = "Russia"
$ws.Range("M1").Value = "Saudi Arabia"
$ws.Range("N1").Value = "USA"

# --- Fix mislabeled row header in A3 ---
$ws.Range("A3").Value = "Affected by the top tax (any variant)"

# --- Row 2: "Supports tax on world top income ... (Any variant)" ---
$ws.Range("B2").Value = 0.562249333952068
$ws.Range("C2").Value = 0.530633111073681
$ws.Range("D2").Value = 0.534615441391404
$ws.Range("E2").Value = 0.532199124020193
$ws.Range("F2").Value = 0.475339526743875
$ws.Range("G2").Value = 0.455066177183276
$ws.Range("H2").Value = 0.481053473314213
$ws.Range("I2").Value = 0.644527456241547
$ws.Range("J2").Value = 0.283838350853901
$ws.Range("K2").Value = 0.278766736483527
$ws.Range("L2").Value = 0.599907829518252
$ws.Range("M2").Value = 0.695880953418338
$ws.Range("N2").Value = 0.38602041410698

# --- Row 3: "Affected by the top tax (any variant)" ---
$ws.Range("B3").Value = 0.0591524213077198
$ws.Range("C3").Value = 0.0407707508788973
$ws.Range("D3").Value = 0.0323072107425518
$ws.Range("E3").Value = 0.0266191432108494
$ws.Range("F3").Value = 0.0508665529273163
$ws.Range("G3").Value = 0.0154979588373725
$ws.Range("H3").Value = 0.0308607490074744
$ws.Range("I3").Value = 0.0702529813060691
$ws.Range("J3").Value = 0.114402739983896
$ws.Range("K3").Value = 0.0300920300421565
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value = 0.213457131326528
$ws.Range("N3").Value = 0.0779223306659385

# --- Row 4: "Supports tax on world top 1% ..." ---
$ws.Range("B4").Value = 0.593568464643459
$ws.Range("C4").Value = 0.623206921829844
$ws.Range("D4").Value = 0.662100661962385
$ws.Range("E4").Value = 0.645129138499509
$ws.Range("F4").Value = 0.652088141940814
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.595169320047074
$ws.Range("I4").Value = 0.581049693515047
$ws.Range("J4").Value = 0.460003610927742
$ws.Range("K4").Value = 0.321840723174092
$ws.Range("L4").Value = 0.599907829518252
$ws.Range("M4").Value = 0.656234738255263
$ws.Range("N4").Value = 0.466960372968151

# --- Row 5: "Affected by the top 1% tax (income > $PPP 120k)" ---
$ws.Range("B5").Value = 0.0324192263079185
$ws.Range("C5").Value = 0.0228809019118487
$ws.Range("D5").Value = 0.0152772811251959
$ws.Range("E5").Value = 0.0158159265407275
$ws.Range("F5").Value = 0.0364460596482764
$ws.Range("G5").Value = 0.00635212756481645
$ws.Range("H5").Value = 0.013703147356521
$ws.Range("I5").Value = 0.0382299828227781
$ws.Range("J5").Value = 0.0699523880243387
$ws.Range("K5").Value = 0.0195161254114759
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = 0.136800251123575
$ws.Range("N5").Value = 0.0397521568871045

# --- Row 6: "Supports tax on world top 3% ..." ---
$ws.Range("B6").Value = 0.564833073218838
$ws.Range("C6").Value = 0.494229631325223
$ws.Range("D6").Value = 0.494509481529061
$ws.Range("E6").Value = 0.482293276239211
$ws.Range("F6").Value = 0.377478972984925
$ws.Range("G6").Value = 0.317062120727541
$ws.Range("H6").Value = 0.446739876773781
$ws.Range("I6").Value = 0.666981999539339
$ws.Range("J6").Value = 0.204001758620865
$ws.Range("K6").Value = 0.259051167955435
$ws.Range("L6").Value = 0.599907829518252
$ws.Range("M6").Value = 0.713245760172143
$ws.Range("N6").Value = 0.359482567475819

# --- Row 7: "Affected by the top 3% tax (income > $PPP 80k)" ---
$ws.Range("B7").Value = 0.085037992635904
$ws.Range("C7").Value = 0.058871310106645
$ws.Range("D7").Value = 0.0497557537097689
$ws.Range("E7").Value = 0.0381284684121215
$ws.Range("F7").Value = 0.0651357948522902
$ws.Range("G7").Value = 0.0243920108908461
$ws.Range("H7").Value = 0.0494954076529742
$ws.Range("I7").Value = 0.0998345665218679
$ws.Range("J7").Value = 0.160672303258471
$ws.Range("K7").Value = 0.0400180177711981
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value = 0.282886812048479
$ws.Range("N7").Value = 0.113725803486659
